# Update transition-probability matrix values on Sheet1 to reflect updated simulation data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2" = 0.1828571428571429
    "C2" = 0.5828571428571429
    "J2" = 0.01333333333333333
    "P2" = 0.1295238095238095
    "S2" = 0.09142857142857143
    "B3" = 0.009554140127388535
    "C3" = 0.01910828025477707
    "J3" = 0.02866242038216561
    "O3" = 0.003184713375796179
    "P3" = 0.7834394904458599
    "S3" = 0.1560509554140127
    "J4" = 0.02531645569620253
    "O4" = 0.01265822784810127
    "P4" = 0.7341772151898734
    "S4" = 0.2278481012658228
    "P5" = 0.5
    "S5" = 0.5
    "B6" = 0.0681265206812652
    "D6" = 0.0072992700729927
    "F6" = 0.08759124087591241
    "J6" = 0.1776155717761557
    "O6" = 0.0218978102189781
    "Q6" = 0.2214111922141119
    "R6" = 0.05596107055961071
    "S6" = 0.3600973236009732
    "B7" = 0.08994708994708994
    "D7" = 0.01587301587301587
    "E7" = 0.005291005291005291
    "F7" = 0.06613756613756613
    "J7" = 0.1111111111111111
    "O7" = 0.02116402116402116
    "Q7" = 0.1931216931216931
    "R7" = 0.07407407407407407
    "S7" = 0.4232804232804233
    "B8" = 0.08801955990220049
    "D8" = 0.02322738386308068
    "E8" = 0.001222493887530562
    "F8" = 0.07946210268948656
    "J8" = 0.09535452322738386
    "O8" = 0.02078239608801956
    "Q8" = 0.176039119804401
    "R8" = 0.09535452322738386
    "S8" = 0.4205378973105134
    "B9" = 0.1040609137055838
    "D9" = 0.01522842639593909
    "F9" = 0.07868020304568528
    "J9" = 0.1116751269035533
    "O9" = 0.01269035532994924
    "Q9" = 0.1700507614213198
    "R9" = 0.1116751269035533
    "S9" = 0.3959390862944163
    "B10" = 0.1053763440860215
    "D10" = 0.02150537634408602
    "E10" = 0.0004301075268817204
    "F10" = 0.0675268817204301
    "J10" = 0.1161290322580645
    "O10" = 0.01806451612903226
    "Q10" = 0.2008602150537634
    "R10" = 0.08516129032258064
    "S10" = 0.3849462365591398
    "G11" = 0.1538461538461539
    "J11" = 0.07871198568872988
    "K11" = 0.2021466905187835
    "L11" = 0.5599284436493739
    "S11" = 0.005366726296958855
    "G12" = 0.764525993883792
    "J12" = 0.1865443425076453
    "K12" = 0.003058103975535168
    "L12" = 0.02140672782874618
    "S12" = 0.02446483180428135
    "G13" = 0.6363636363636364
    "J13" = 0.3068181818181818
    "S13" = 0.05681818181818182
    "F14" = 0.2
    "G14" = 0.8
    "F15" = 0.01168224299065421
    "H15" = 0.1495327102803738
    "I15" = 0.05841121495327103
    "J15" = 0.3995327102803738
    "K15" = 0.06074766355140187
    "M15" = 0.01635514018691589
    "O15" = 0.06542056074766354
    "S15" = 0.2383177570093458
    "F16" = 0.01400560224089636
    "H16" = 0.1764705882352941
    "I16" = 0.09523809523809523
    "J16" = 0.42296918767507
    "K16" = 0.1064425770308123
    "M16" = 0.01680672268907563
    "O16" = 0.05042016806722689
    "S16" = 0.1176470588235294
    "F17" = 0.01311084624553039
    "H17" = 0.1632896305125149
    "I17" = 0.09535160905840286
    "J17" = 0.4314660309892729
    "K17" = 0.09773539928486293
    "M17" = 0.02264600715137068
    "N17" = 0.003575685339690107
    "O17" = 0.05721096543504171
    "S17" = 0.1156138259833135
    "F18" = 0.002702702702702703
    "H18" = 0.1540540540540541
    "I18" = 0.0918918918918919
    "J18" = 0.4405405405405405
    "K18" = 0.1027027027027027
    "M18" = 0.01621621621621622
    "O18" = 0.06216216216216217
    "S18" = 0.1297297297297297
    "F19" = 0.01278350515463918
    "H19" = 0.2078350515463918
    "I19" = 0.09072164948453608
    "J19" = 0.3542268041237113
    "K19" = 0.1063917525773196
    "M19" = 0.02268041237113402
    "N19" = 0.0008247422680412372
    "O19" = 0.07134020618556701
    "S19" = 0.1331958762886598
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
